# Updates cryptos list values (Price / Volume(1h)) to match the latest
# scrape, and fixes the Stacks / WhiteBITCoin row ordering (rows 45-46).
# All target cells are forced to Text format before assignment so that
# numeric-looking strings (e.g. "5.00", "1.944.06") are preserved exactly
# as text, matching the original inline-string cell content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '54.540.10'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -6.90%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.428.19'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -10.91%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '465.83'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -7.53%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.44'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -5.19%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.20%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.493'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -6.87%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.450.75'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -10.46%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0953'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -8.85%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.32'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -12.02%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.315'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -8.95%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.122'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.89%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.857.12'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -10.84%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '54.417.46'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -7.34%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000133'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.52%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.81'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -8.49%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.450.36'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -10.24%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.19'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -11.83%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '310.32'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -9.59%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.49'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -13.42%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.997'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.11%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.92%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.38'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -13.89%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '56.25'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -10.56%  '

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.80%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.385'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -9.50%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.156'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -9.47%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.520.63'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -11.51%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.13'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -5.07%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.997'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.15%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0₃0719'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -13.13%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '146.86'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.44%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.81'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -6.93%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -10.31%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.00'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -7.64%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.55'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -15.16%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.06'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -6.21%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.809'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -14.64%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.995'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.12%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '32.95'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -7.93%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.597'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.97%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0524'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -6.12%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.25'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -8.26%  '

$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.24'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -10.62%  '

$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.09'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.68%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.944.06'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -11.18%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0885'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.02%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0218'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.91%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.24'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -10.59%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.70'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -12.03%  '
